$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BeiJing
$ws.Range("C2").Value = 10672
$ws.Range("D2").Value = 77
$ws.Range("E2").Value = 896

# Row 3 - ChengDu
$ws.Range("C3").Value = 1455
$ws.Range("E3").Value = 140

# Row 4 - ShangHai
$ws.Range("C4").Value = 7030

# Row 5 - GuangZhou
$ws.Range("C5").Value = 1966
$ws.Range("E5").Value = 133

# Row 6 - ShenYang
$ws.Range("C6").Value = 915
$ws.Range("E6").Value = 101

# Row 7 - Europe
$ws.Range("C7").Value = 354
